$d = $word.ActiveDocument
$t = $d.Tables(1)

$replacements = @(
    @(1, 1, "88÷5=17, 3"),
    @(1, 2, "27÷2=13, 1"),
    @(1, 3, "37÷9=4, 1"),
    @(1, 4, "80÷4=20, 0"),
    @(1, 5, "64÷4=16, 0"),
    @(5, 1, "21÷7=3, 0"),
    @(5, 2, "18÷4=4, 2"),
    @(5, 3, "75÷5=15, 0"),
    @(5, 4, "63÷6=10, 3"),
    @(5, 5, "79÷2=39, 1"),
    @(9, 1, "46÷5=9, 1"),
    @(9, 2, "29÷5=5, 4"),
    @(9, 3, "29÷5=5, 4"),
    @(9, 4, "42÷6=7, 0"),
    @(9, 5, "96÷9=10, 6"),
    @(13, 1, "42÷5=8, 2"),
    @(13, 2, "58÷8=7, 2"),
    @(13, 3, "94÷2=47, 0"),
    @(13, 4, "89÷9=9, 8"),
    @(13, 5, "50÷6=8, 2"),
    @(17, 1, "39÷8=4, 7"),
    @(17, 2, "53÷7=7, 4"),
    @(17, 3, "26÷6=4, 2"),
    @(17, 4, "91÷4=22, 3"),
    @(17, 5, "50÷7=7, 1")
)

foreach ($pair in $replacements) {
    $row = $pair[0]
    $col = $pair[1]
    $new = $pair[2]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $new
}
